$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: rows 2-14 (introduces SECBM, SECBL shared strings) ---
$ws.Cells.Item(2,1).Value = "AB_ST_EXISTING"
$ws.Cells.Item(2,2).Value = "Steam Turbine Using Agricultural By-Products (EIA 860 Nomenclature)"
$ws.Cells.Item(2,3).Value = "SECBM"
$ws.Cells.Item(3,1).Value = "BIT_ST_EXISTING"
$ws.Cells.Item(3,2).Value = "Steam Turbine Using Bituminous Coal (EIA 860 Nomenclature)"
$ws.Cells.Item(3,3).Value = "SECBM"
$ws.Cells.Item(4,1).Value = "BLQ_ST_EXISTING"
$ws.Cells.Item(4,2).Value = "Steam Turbine Using Black Liquor (EIA 860 Nomenclature)"
$ws.Cells.Item(4,3).Value = "SECBM"
$ws.Cells.Item(5,1).Value = "DFO_CC_EXISTING"
$ws.Cells.Item(5,2).Value = "Combined Cycle Combustion Turbine Using Petroleum (EIA 860 Nomenclature)"
$ws.Cells.Item(5,3).Value = "SECBM"
$ws.Cells.Item(6,1).Value = "DFO_GT_EXISTING"
$ws.Cells.Item(6,2).Value = "Combustion Turbine Using Petroleum (EIA 860 Nomenclature)"
$ws.Cells.Item(6,3).Value = "SECBM"
$ws.Cells.Item(7,1).Value = "DFO_IC_EXISTING"
$ws.Cells.Item(7,2).Value = "Internal Combustion Engine Using Petroleum (EIA 860 Nomenclature)"
$ws.Cells.Item(7,3).Value = "SECBM"
$ws.Cells.Item(8,1).Value = "LFG_GT_EXISTING"
$ws.Cells.Item(8,2).Value = "Combustion Turbine Using Landfill Gas (EIA 860 Nomenclature)"
$ws.Cells.Item(8,3).Value = "SECBM"
$ws.Cells.Item(9,1).Value = "LFG_IC_EXISTING"
$ws.Cells.Item(9,2).Value = "Internal Combustion Engine Using Landfill Gas (EIA 860 Nomenclature)"
$ws.Cells.Item(9,3).Value = "SECBM"
$ws.Cells.Item(10,1).Value = "MWH_BA1H_EXISTING"
$ws.Cells.Item(10,2).Value = "Battery  Storage- 1h  (EIA 860 Nomenclature)"
$ws.Cells.Item(10,3).Value = "SECBL"
$ws.Cells.Item(11,1).Value = "MWH_BA2H_EXISTING"
$ws.Cells.Item(11,2).Value = "Battery  Storage- 2h  (EIA 860 Nomenclature)"
$ws.Cells.Item(11,3).Value = "SECBL"
$ws.Cells.Item(12,1).Value = "NG_CC_EXISTING"
$ws.Cells.Item(12,2).Value = "Combined Cycle Combustion Turbine Using Natural Gas (EIA 860 Nomenclature)"
$ws.Cells.Item(12,3).Value = "SECBM"
$ws.Cells.Item(13,1).Value = "NG_GT_EXISTING"
$ws.Cells.Item(13,2).Value = "Combustion Turbine Using Natural Gas (EIA 860 Nomenclature)"
$ws.Cells.Item(13,3).Value = "SECBM"
$ws.Cells.Item(14,1).Value = "NG_ST_EXISTING"
$ws.Cells.Item(14,2).Value = "Steam Turbine Using Natural Gas (EIA 860 Nomenclature)"
$ws.Cells.Item(14,3).Value = "SECBM"

# --- Phase 2: fix Tech Description for BIOMASS_NEW (row 28) ahead of the rest ---
$ws.Cells.Item(28,2).Value = "Generation From Biomass No Carbon Capture"

# --- Phase 3: remaining rows 15-27 and 29-56 (in order) ---
$ws.Cells.Item(15,1).Value = "NUC_ST_EXISTING"
$ws.Cells.Item(15,2).Value = "Nuclear Turbine (EIA 860 Nomenclature)"
$ws.Cells.Item(15,3).Value = "CECBM"
$ws.Cells.Item(16,1).Value = "OBG_IC_EXISTING"
$ws.Cells.Item(16,2).Value = "Internal Combustion Engine Using Other Biomass Gas (EIA 860 Nomenclature)"
$ws.Cells.Item(16,3).Value = "SECBM"
$ws.Cells.Item(17,1).Value = "SUN_PV_EXISTING"
$ws.Cells.Item(17,2).Value = "Solar Photovoltaic - Utility (EIA 860 Nomenclature)"
$ws.Cells.Item(17,3).Value = "solar_Kabre"
$ws.Cells.Item(18,1).Value = "WAT_HY_EXISTING"
$ws.Cells.Item(18,2).Value = "Conventional Hydroelectric (EIA 860 Nomenclature)"
$ws.Cells.Item(18,3).Value = "CECBM"
$ws.Cells.Item(19,1).Value = "WAT_PS_EXISTING"
$ws.Cells.Item(19,2).Value = "Hydroelectric Pumped Storage (EIA 860 Nomenclature)"
$ws.Cells.Item(19,3).Value = "CECBM"
$ws.Cells.Item(20,1).Value = "WDS_ST_EXISTING"
$ws.Cells.Item(20,2).Value = "Steam Turbine Using Wood Waste (EIA 860 Nomenclature)"
$ws.Cells.Item(20,3).Value = "SECBM"
$ws.Cells.Item(21,1).Value = "WH_ST_EXISTING"
$ws.Cells.Item(21,2).Value = "Steam Turbine Using Waste Heat (EIA 860 Nomenclature)"
$ws.Cells.Item(21,3).Value = "SECBM"
$ws.Cells.Item(22,1).Value = "WND_WT_EXISTING"
$ws.Cells.Item(22,2).Value = "Onshore Wind Turbine (EIA 860 Nomenclature)"
$ws.Cells.Item(22,3).Value = "wind_AVG_CAMPO_MIGUEL"
$ws.Cells.Item(23,1).Value = "BATT_2H_NEW"
$ws.Cells.Item(23,2).Value = "Battery Storage 2h (NREL ATB 2023 Technology)"
$ws.Cells.Item(23,3).Value = "SECBL"
$ws.Cells.Item(24,1).Value = "BATT_4H_NEW"
$ws.Cells.Item(24,2).Value = "Battery Storage 4h (NREL ATB 2023 Technology)"
$ws.Cells.Item(24,3).Value = "SECBL"
$ws.Cells.Item(25,1).Value = "BATT_6H_NEW"
$ws.Cells.Item(25,2).Value = "Battery Storage 6h (NREL ATB 2023 Technology)"
$ws.Cells.Item(25,3).Value = "SECBL"
$ws.Cells.Item(26,1).Value = "BATT_8H_NEW"
$ws.Cells.Item(26,2).Value = "Battery Storage 8h (NREL ATB 2023 Technology)"
$ws.Cells.Item(26,3).Value = "SECBL"
$ws.Cells.Item(27,1).Value = "BIOMASS_CC90_NEW"
$ws.Cells.Item(27,2).Value = "Generation From Biomass With 90% Carbon Capture (Technology from NREL ReEDS model  Using BECC-mod)"
$ws.Cells.Item(27,3).Value = "SECBM"
$ws.Cells.Item(29,1).Value = "CO2_STORAGE"
$ws.Cells.Item(29,2).Value = "CO2 Storage"
$ws.Cells.Item(29,3).Value = "NotAffected"
$ws.Cells.Item(30,1).Value = "COAL_95CC_NEW"
$ws.Cells.Item(30,2).Value = "Generation From Coal With 95% Carbon Capture (NREL ATB 2023 Technology)"
$ws.Cells.Item(30,3).Value = "SECBM"
$ws.Cells.Item(31,1).Value = "COAL_99CC_NEW"
$ws.Cells.Item(31,2).Value = "Generation From Coal With 99% Carbon Capture (NREL ATB 2023 Technology)"
$ws.Cells.Item(31,3).Value = "SECBM"
$ws.Cells.Item(32,1).Value = "COAL_NEW"
$ws.Cells.Item(32,2).Value = "Generation From Coal (NREL ATB 2023 Technology)"
$ws.Cells.Item(32,3).Value = "SECBM"
$ws.Cells.Item(33,1).Value = "FT_BIOMASS"
$ws.Cells.Item(33,2).Value = "Fuel for Generation Technologies that Use Biomass"
$ws.Cells.Item(33,3).Value = "NotAffected"
$ws.Cells.Item(34,1).Value = "FT_COAL"
$ws.Cells.Item(34,2).Value = "Fuel for Generation Technologies that Use Coal"
$ws.Cells.Item(34,3).Value = "NotAffected"
$ws.Cells.Item(35,1).Value = "FT_NG"
$ws.Cells.Item(35,2).Value = "Fuel for Generation Technologies that Use Natural Gas"
$ws.Cells.Item(35,3).Value = "NotAffected"
$ws.Cells.Item(36,1).Value = "FT_NUCLEAR"
$ws.Cells.Item(36,2).Value = "Fuel for Nuclear Generation Technologies"
$ws.Cells.Item(36,3).Value = "NotAffected"
$ws.Cells.Item(37,1).Value = "FT_PETROLEUM"
$ws.Cells.Item(37,2).Value = "Fuel for Generation Technologies that Use Petroleum"
$ws.Cells.Item(37,3).Value = "NotAffected"
$ws.Cells.Item(38,1).Value = "NG_F-FRAME_CC_95CC_NEW"
$ws.Cells.Item(38,2).Value = "Combined Cycle Natural Gas Turbine F-Frame With 95 % of Carbon Capture (NREL ATB 2023 Technology)"
$ws.Cells.Item(38,3).Value = "SECBM"
$ws.Cells.Item(39,1).Value = "NG_F-FRAME_CC_97CC_NEW"
$ws.Cells.Item(39,2).Value = "Combined Cycle Natural Gas Turbine F-Frame With 97 % of Carbon Capture (NREL ATB 2023 Technology)"
$ws.Cells.Item(39,3).Value = "SECBM"
$ws.Cells.Item(40,1).Value = "NG_F-FRAME_CC_NEW"
$ws.Cells.Item(40,2).Value = "Combined Cycle Natural Gas Turbine F-Frame (NREL ATB 2023 Technology)"
$ws.Cells.Item(40,3).Value = "SECBM"
$ws.Cells.Item(41,1).Value = "NG_F-FRAME_CT_NEW"
$ws.Cells.Item(41,2).Value = "Natural Gas Combustion Turbine F-Frame - Simple Cycle (NREL ATB 2023 Technology)"
$ws.Cells.Item(41,3).Value = "SECBM"
$ws.Cells.Item(42,1).Value = "NG_H-FRAME_CC_95CC_NEW"
$ws.Cells.Item(42,2).Value = "Combined Cycle Natural Gas Turbine H-Frame With 95 % of Carbon Capture (NREL ATB 2023 Technology)"
$ws.Cells.Item(42,3).Value = "SECBM"
$ws.Cells.Item(43,1).Value = "NG_H-FRAME_CC_97CC_NEW"
$ws.Cells.Item(43,2).Value = "Combined Cycle Natural Gas Turbine H-Frame With 97 % of Carbon Capture (NREL ATB 2023 Technology)"
$ws.Cells.Item(43,3).Value = "SECBM"
$ws.Cells.Item(44,1).Value = "NG_H-FRAME_CC_NEW"
$ws.Cells.Item(44,2).Value = "Combined Cycle Natural Gas Turbine H-Frame (NREL ATB 2023 Technology)"
$ws.Cells.Item(44,3).Value = "SECBM"
$ws.Cells.Item(45,1).Value = "NUCLEAR-AP1000_NEW"
$ws.Cells.Item(45,2).Value = "Nuclear Generation Using AP1000 PWR (NREL ATB 2023 Technology)"
$ws.Cells.Item(45,3).Value = "CECBM"
$ws.Cells.Item(46,1).Value = "NUCLEAR-SMR_NEW"
$ws.Cells.Item(46,2).Value = "Small Modular Nuclear Reactor (NREL ATB 2023 Technology)"
$ws.Cells.Item(46,3).Value = "CECBM"
$ws.Cells.Item(47,1).Value = "PV-COMMERCIAL_NEW"
$ws.Cells.Item(47,2).Value = "Commercial Solar PV (NREL ATB 2023 Technology)"
$ws.Cells.Item(47,3).Value = "solar_Kabre"
$ws.Cells.Item(48,1).Value = "PV-RESIDENTIAL_NEW"
$ws.Cells.Item(48,2).Value = "Residential Solar PV (NREL ATB 2023 Technology)"
$ws.Cells.Item(48,3).Value = "solar_Kabre"
$ws.Cells.Item(49,1).Value = "PV-UTILITY_NEW"
$ws.Cells.Item(49,2).Value = "Utility Solar PV (NREL ATB 2023 Technology)"
$ws.Cells.Item(49,3).Value = "solar_Kabre"
$ws.Cells.Item(50,1).Value = "WAT_HY_NEW"
$ws.Cells.Item(50,2).Value = "Conventional Hydroelectric (NREL ATB 2023 Technology)"
$ws.Cells.Item(50,3).Value = "CECBM"
$ws.Cells.Item(51,1).Value = "WAT_PS_NEW"
$ws.Cells.Item(51,2).Value = "Hydroelectric Pumped Storage (NREL ATB 2023 Technology)"
$ws.Cells.Item(51,3).Value = "CECBM"
$ws.Cells.Item(52,1).Value = "WIND-LAND-C8_NEW"
$ws.Cells.Item(52,2).Value = "Onshore Wind Turbine Class 8 From NREL ATB 2023  (NREL ATB 2023 Technology)"
$ws.Cells.Item(52,3).Value = "wind_AVG_CAMPO_MIGUEL"
$ws.Cells.Item(53,1).Value = "WIND-OFFSHORE-C6_NEW"
$ws.Cells.Item(53,2).Value = "Offshore Wind Turbine Class 6 From NREL ATB 2023  (NREL ATB 2023 Technology)"
$ws.Cells.Item(53,3).Value = "wind_AVG_CAMPO_MIGUEL"
$ws.Cells.Item(54,1).Value = "DISTRIBUTION"
$ws.Cells.Item(54,2).Value = "Energy Distribution"
$ws.Cells.Item(54,3).Value = "NotAffected"
$ws.Cells.Item(55,1).Value = "TRANSMISSION_INTERREGIONAL"
$ws.Cells.Item(55,2).Value = "Transmission Between Different Regions"
$ws.Cells.Item(55,3).Value = "NotAffected"
$ws.Cells.Item(56,1).Value = "TRANSMISSION_REGIONAL"
$ws.Cells.Item(56,2).Value = "Transmission In the Same Region"
$ws.Cells.Item(56,3).Value = "NotAffected"

# --- Phase 4: finish row 28 (A and C columns) ---
$ws.Cells.Item(28,1).Value = "BIOMASS_NEW"
$ws.Cells.Item(28,3).Value = "SECBM"

# Apply red fill to the moved Distribution/Transmission rows (A:B only)
$ws.Range("A54:B56").Interior.Color = 255

# Update selection to mirror final interactive state
$ws.Range("A54:XFD56").Select() | Out-Null

# Set page orientation to portrait (page setup)
$ws.PageSetup.Orientation = 1
